# 2021-03-17_Case List에서 일치 Case 찾기
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo / formatting of the Elm 3DS opposite-party name (was "ELM3DS")
$ws.Range("B3").Value = "Elm 3DS"

# Add a new "Note" column (C) that records the matched Case type for each row
$ws.Range("C1").Value = "Note"
$ws.Range("C3").Value = "Patent"
$ws.Range("C4").Value = "Patent"

# Resize column B to fit the new, wider contents and refresh the active selection
$ws.Columns("B").ColumnWidth = 17.5
$ws.Range("F10").Select() | Out-Null
